$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.428715350796693
$ws.Range("D2").Value = [double]"6.505707639270497E-08"
$ws.Range("C3").Value = 7.588202386008843
$ws.Range("D3").Value = [double]"5.351274978693255E-14"
$ws.Range("C4").Value = 14.95605546505332
$ws.Range("C5").Value = 11.38019880374943
$ws.Range("C6").Value = 12.25143824861036
$ws.Range("C7").Value = 11.19106500110048
$ws.Range("C8").Value = 13.79811665608808
$ws.Range("C9").Value = 4.542554067930547
$ws.Range("D9").Value = [double]"5.954015320330797E-06"
$ws.Range("C10").Value = 8.026149282465457
$ws.Range("D10").Value = [double]"1.998401444325282E-15"
$ws.Range("C11").Value = 12.26251647594378
$ws.Range("C12").Value = 10.15780203888068
$ws.Range("C13").Value = 10.26720719954728
$ws.Range("C14").Value = 9.45488020867294
$ws.Range("C15").Value = 11.3616561052871
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = 5.687090854113968
$ws.Range("D16").Value = [double]"1.521301817852816E-08"
$ws.Range("C17").Value = 7.533170029143039
$ws.Range("D17").Value = [double]"8.060219158778636E-14"
$ws.Range("C18").Value = 15.00383663467488
$ws.Range("C19").Value = 11.29321781418928
$ws.Range("C20").Value = 12.17023307139602
$ws.Range("C21").Value = 11.19092885148498
$ws.Range("C22").Value = 13.79717325390092
$ws.Range("C23").Value = 3.998669556398251
$ws.Range("D23").Value = [double]"6.646301049251235E-05"
$ws.Range("C24").Value = 8.2590755957063
$ws.Range("D24").Value = [double]"2.220446049250313E-16"
$ws.Range("C25").Value = 13.83915927442726
$ws.Range("C26").Value = 9.902007122139366
$ws.Range("C27").Value = 11.2694593108779
$ws.Range("C28").Value = 10.99490916883185
$ws.Range("C29").Value = 13.62872290884902
$ws.Range("C30").Value = 11.7288793133697
$ws.Range("C31").Value = 9.679459355247777
$ws.Range("C32").Value = 12.73570378805782
$ws.Range("C33").Value = 11.85234556241467
$ws.Range("C34").Value = 10.39556310996792
$ws.Range("C35").Value = 8.729830154733497
$ws.Range("C36").Value = 11.40908217091391
$ws.Range("C37").Value = 11.21394411130164
$ws.Range("C38").Value = 8.600704006392009
$ws.Range("D38").Value = 0
$ws.Range("C39").Value = 11.21083720105544
$ws.Range("C40").Value = 12.22947840261592
$ws.Range("C41").Value = 11.43262446024979
$ws.Range("C42").Value = 9.102135207974369
$ws.Range("C43").Value = 13.04831366252636
$ws.Range("C44").Value = 10.22083734255848
$ws.Range("C45").Value = 6.680937001884318
$ws.Range("D45").Value = [double]"3.219224886663596E-11"
$ws.Range("C46").Value = 10.66849421365548
$ws.Range("C47").Value = 9.024525614698298
$ws.Range("D47").Value = 0
$ws.Range("C48").Value = 9.13830197073249
$ws.Range("C49").Value = 8.605860809629439
$ws.Range("C50").Value = 10.15720911016083
$ws.Range("D50").Value = 0
$ws.Range("C51").Value = 5.607504312747027
$ws.Range("D51").Value = [double]"2.395753306849713E-08"
$ws.Range("C52").Value = 7.826943042716973
$ws.Range("D52").Value = [double]"8.659739592076221E-15"
$ws.Range("C53").Value = 13.90189772327733
$ws.Range("C54").Value = 11.15429726278328
$ws.Range("C55").Value = 11.66045329762489
$ws.Range("C56").Value = 10.56870668144853
$ws.Range("C57").Value = 12.87729699401558
$ws.Range("C58").Value = 0.3211630137570719
$ws.Range("D58").Value = 0.74812683126579
$ws.Range("C59").Value = -1.866949510347526
$ws.Range("D59").Value = 0.06208280264152122
$ws.Range("C60").Value = -1.498646340472044
$ws.Range("D60").Value = 0.1341533866968256
$ws.Range("C61").Value = 1.70259819380042
$ws.Range("D61").Value = 0.08882838250278091
$ws.Range("E61").Value = "No"
$ws.Range("C62").Value = 4.162202535518333
$ws.Range("D62").Value = [double]"3.311586657472887E-05"
$ws.Range("C63").Value = 12.94186943227849
$ws.Range("C64").Value = 10.36664986293461
